# Catalog.xlsx edit: add the "AddPriceAgrmnt_FavFolder" worksheet after the
# last existing sheet (AddNonPriceAgr_GlobalCatalog), populate it with the
# Role/Location/selectUOM header row + sample data row, size its columns,
# and leave the selection on C5 (matching the freshly-authored sheet in
# the target workbook). This naturally clears tabSelected on the
# previously-last sheet and shifts the workbook's active tab index.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "AddPriceAgrmnt_FavFolder"

$ws.Range("A1").Value = "Role"
$ws.Range("B1").Value = "Location                      "
$ws.Range("C1").Value = "selectUOM "

$ws.Range("A2").Value = "REQUESTOR"
$ws.Range("B2").Value = "XEEVA -MJ"
$ws.Range("C2").Value = "CU-CUBIC"

$ws.Columns.Item(2).ColumnWidth = 9.7109375
$ws.Columns.Item(3).ColumnWidth = 11.140625

$ws.Range("C5").Select() | Out-Null
